$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.517.77"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "3.556.74"
$ws.Range("E3").Value = "  +3.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.41"
$ws.Range("E5").Value = "  +1.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.90"
$ws.Range("E6").Value = "  +3.45%  "

$ws.Range("D7").Value = "3.557.32"
$ws.Range("E7").Value = "  +3.28%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.126"
$ws.Range("E10").Value = "  +2.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.07"
$ws.Range("E11").Value = "  -4.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.396"
$ws.Range("E12").Value = "  +4.34%  "

$ws.Range("D13").Value = "4.160.30"
$ws.Range("E13").Value = "  +3.31%  "

$ws.Range("E14").Value = "  +3.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.27"
$ws.Range("E15").Value = "  +2.14%  "

$ws.Range("D16").Value = "3.566.98"
$ws.Range("E16").Value = "  +3.92%  "

$ws.Range("D18").Value = "65.539.06"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.39"
$ws.Range("E19").Value = "  +4.83%  "

$ws.Range("E20").Value = "  +1.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.32"
$ws.Range("E21").Value = "  +4.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "395.55"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("E23").Value = "  +4.73%  "

$ws.Range("D24").Value = "3.700.43"
$ws.Range("E24").Value = "  +3.16%  "

$ws.Range("E25").Value = "  +1.12%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  +10.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.89"
$ws.Range("E28").Value = "  +9.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("E30").Value = "  +1.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.43"
$ws.Range("E31").Value = "  +2.54%  "

$ws.Range("D32").Value = "3.570.97"
$ws.Range("E32").Value = "  +3.49%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.147"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.84"
$ws.Range("E35").Value = "  +3.60%  "

$ws.Range("E36").Value = "  +6.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.07"
$ws.Range("E37").Value = "  +1.97%  "

$ws.Range("E38").Value = "  +3.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "167.54"
$ws.Range("E39").Value = "  -2.56%  "

$ws.Range("E40").Value = "  +4.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0811"
$ws.Range("E41").Value = "  +4.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.834"
$ws.Range("E42").Value = "  +1.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.69"
$ws.Range("E43").Value = "  +16.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.11"
$ws.Range("E44").Value = "  -0.81%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.46"
$ws.Range("E46").Value = "  +0.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.70"
$ws.Range("E47").Value = "  +5.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.21"
$ws.Range("E48").Value = "  +8.94%  "

$ws.Range("D49").Value = "2.486.08"
$ws.Range("E49").Value = "  +12.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.84"
$ws.Range("E50").Value = "  +3.72%  "

# Row 51: coin changed from LidoDAOToken to dogwifhat
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.10"
$ws.Range("E51").Value = "  -0.04%  "
